$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet 1: CHI TIẾT DOANH THU
#   Insert two new revenue rows (620, 621) before the "Tổng" row,
#   pushing the totals row from row4 -> row6 and updating its totals.
# ======================================================================
$ws1 = $wb.Worksheets.Item("CHI TIẾT DOANH THU")
$ws1.Rows("4:5").Insert()

$ws1.Range("A4:B5").NumberFormat = "@"

$ws1.Range("A4").Value = "08-03-2024"
$ws1.Range("B4").Value = "HD-LUXURY"
$ws1.Range("C4").Value = 620
$ws1.Range("D4").Value = "CẦN THƠ"
$ws1.Range("E4").Value = "Nâng mũi"
$ws1.Range("F4").Value = "Trần Thị Thanh Nhàn"
$ws1.Range("G4").Value = "Cá nhân"
$ws1.Range("H4").Value = "Phạm Thanh Hoàng"
$ws1.Range("I4").Value = 28000000
$ws1.Range("L4").Value = 28000000
$ws1.Range("M4").Value = "Phạm Thanh Hoàng"
$ws1.Range("O4").Value = 28000000
$ws1.Range("P4").Value = 0
$ws1.Range("Q4").Value = 28000000
$ws1.Range("R4").Value = 0
$ws1.Range("S4").Value = "Lâm Hoàng Phú"
$ws1.Range("U4").Value = 100000
$ws1.Range("V4").Value = 50000

$ws1.Range("A5").Value = "08-03-2024"
$ws1.Range("B5").Value = "HD-LUXURY"
$ws1.Range("C5").Value = 621
$ws1.Range("D5").Value = "CẦN THƠ"
$ws1.Range("E5").Value = "Nâng mũi"
$ws1.Range("F5").Value = "Trần Thị Ngọc Dung"
$ws1.Range("G5").Value = "Cá nhân"
$ws1.Range("H5").Value = "Lâm Thị Mỹ Hằng"
$ws1.Range("I5").Value = 9000000
$ws1.Range("L5").Value = 9000000
$ws1.Range("M5").Value = "Phạm Thanh Hoàng"
$ws1.Range("O5").Value = 9000000
$ws1.Range("P5").Value = 0
$ws1.Range("Q5").Value = 9000000
$ws1.Range("R5").Value = 0
$ws1.Range("S5").Value = "Lâm Hoàng Phú"
$ws1.Range("U5").Value = 100000
$ws1.Range("V5").Value = 50000

# Update the "Tổng" row (now row 6)
$ws1.Range("C6").Value = 4
$ws1.Range("I6").Value = 39100000
$ws1.Range("L6").Value = 45100000
$ws1.Range("O6").Value = 45100000
$ws1.Range("Q6").Value = 45100000
$ws1.Range("U6").Value = 300000
$ws1.Range("V6").Value = 100000

# ======================================================================
# Sheet 2: CHI TIẾT VỀ THU NỢ
#   Insert a new "Ngày thực hiện" column before "Lượng thu", add three
#   debt-collection rows, and push + update the "Tổng" row to row5.
# ======================================================================
$ws2 = $wb.Worksheets.Item("CHI TIẾT VỀ THU NỢ")
$ws2.Columns("F:F").Insert()
$ws2.Range("F1").Value = "Ngày thực hiện"

$ws2.Rows("2:4").Insert()

$ws2.Range("C2:C4").NumberFormat = "@"
$ws2.Range("F2:F4").NumberFormat = "@"

$ws2.Range("A2").Value = "TN"
$ws2.Range("B2").Value = 171
$ws2.Range("C2").Value = "08-03-2024"
$ws2.Range("D2").Value = "CẦN THƠ"
$ws2.Range("E2").Value = "HD-LUXURY-190"
$ws2.Range("F2").Value = "2024-04-16"
$ws2.Range("G2").Value = 1000000

$ws2.Range("A3").Value = "TN"
$ws2.Range("B3").Value = 172
$ws2.Range("C3").Value = "08-03-2024"
$ws2.Range("D3").Value = "CẦN THƠ"
$ws2.Range("E3").Value = "HD-LUXURY-191"
$ws2.Range("F3").Value = "2024-04-16"
$ws2.Range("G3").Value = 4000000

$ws2.Range("A4").Value = "TN"
$ws2.Range("B4").Value = 173
$ws2.Range("C4").Value = "08-03-2024"
$ws2.Range("D4").Value = "CẦN THƠ"
$ws2.Range("E4").Value = "HD-LUXURY-587"
$ws2.Range("F4").Value = "2024-07-24"
$ws2.Range("G4").Value = 10000000

# Update the "Tổng" row (now row 5)
$ws2.Range("B5").Value = 3
$ws2.Range("F5").Value = ""
$ws2.Range("G5").Value = 15000000

# ======================================================================
# Sheet 3: CHI TIẾT CHI TIÊU
#   Add two new expense rows (750 Ứng Lương, 751 Chi Phí Sinh Hoạt Tại
#   Cơ Sở) before the "Tổng" row, pushing it to row7 with new totals.
# ======================================================================
$ws3 = $wb.Worksheets.Item("CHI TIẾT CHI TIÊU")
$ws3.Rows("5:6").Insert()

$ws3.Range("C5:C6").NumberFormat = "@"

$ws3.Range("A5").Value = "CT"
$ws3.Range("B5").Value = 750
$ws3.Range("C5").Value = "08-03-2024"
$ws3.Range("D5").Value = "CẦN THƠ"
$ws3.Range("E5").Value = "Ứng Lương"
$ws3.Range("F5").Value = 1000000

$ws3.Range("A6").Value = "CT"
$ws3.Range("B6").Value = 751
$ws3.Range("C6").Value = "08-03-2024"
$ws3.Range("D6").Value = "CẦN THƠ"
$ws3.Range("E6").Value = "Chi Phí Sinh Hoạt Tại Cơ Sở"
$ws3.Range("F6").Value = 2230000

# Update the "Tổng" row (now row 7)
$ws3.Range("B7").Value = 5
$ws3.Range("F7").Value = 10593000

# ======================================================================
# Sheet 4: DOANH SỐ CÁ NHÂN  (value-only updates, no new rows)
# ======================================================================
$ws4 = $wb.Worksheets.Item("DOANH SỐ CÁ NHÂN")
$ws4.Range("F3").Value = 3
$ws4.Range("G3").Value = 250000
$ws4.Range("B4").Value = 9000000
$ws4.Range("J4").Value = 137600000
$ws4.Range("J6").Value = 28000000
$ws4.Range("B9").Value = 28000000
$ws4.Range("D9").Value = 37000000
$ws4.Range("B12").Value = 39100000
$ws4.Range("D12").Value = 45100000
$ws4.Range("F12").Value = 3
$ws4.Range("G12").Value = 250000
$ws4.Range("J12").Value = 247600000

# ======================================================================
# Sheet 5: CHI TIÊU TỔNG HỢP
#   Re-sorted category breakdown with two new categories inserted.
# ======================================================================
$ws5 = $wb.Worksheets.Item("CHI TIÊU TỔNG HỢP")
$ws5.Rows("2:2").Insert()
$ws5.Rows("5:5").Insert()

$ws5.Range("A2").Value = "Chi Phí Sinh Hoạt Tại Cơ Sở"
$ws5.Range("B2").Value = 2230000
$ws5.Range("A3").Value = "Chi Phí Vận Hành"
$ws5.Range("B3").Value = 220000
$ws5.Range("A4").Value = "Trang thiết bị Y Tế"
$ws5.Range("B4").Value = 2143000
$ws5.Range("A5").Value = "Ứng Lương"
$ws5.Range("B5").Value = 1000000
$ws5.Range("A6").Value = "Blank"
$ws5.Range("B6").Value = 5000000
$ws5.Range("A7").Value = "Tổng cộng"
$ws5.Range("B7").Value = 10593000

# ======================================================================
# Sheet 6: LŨY KẾ NGÀY
#   Convert old "Tổng" row (row4) into an "08-03-2024" data row, and
#   append a fresh "Tổng" row at row5.
# ======================================================================
$ws6 = $wb.Worksheets.Item("LŨY KẾ NGÀY")
$ws6.Rows("5:5").Insert()

$ws6.Range("A4").NumberFormat = "@"
$ws6.Range("A4").Value = "08-03-2024"
$ws6.Range("B4").Value = 37000000
$ws6.Range("C4").Value = 37000000
$ws6.Range("D4").Value = 2
$ws6.Range("E4").Value = 15000000
$ws6.Range("F4").Value = 3230000
$ws6.Range("G4").Value = 48770000

$ws6.Range("A5").Value = "Tổng"
$ws6.Range("B5").Value = 45100000
$ws6.Range("C5").Value = 45100000
$ws6.Range("D5").Value = 4
$ws6.Range("E5").Value = 15000000
$ws6.Range("F5").Value = 10593000
$ws6.Range("G5").Value = 49507000

# ======================================================================
# Sheet 7: QUỸ LƯƠNG  (value-only updates, no new rows)
# ======================================================================
$ws7 = $wb.Worksheets.Item("QUỸ LƯƠNG")
$ws7.Range("C4").Value = 482142.8571428572
$ws7.Range("C10").Value = 2378095.238095238
$ws7.Range("C11").Value = 7771071.428571429
$ws7.Range("C12").Value = 2141428.571428571
$ws7.Range("C23").Value = 16225023.80952381

# ======================================================================
# Sheet 8: LỢI NHUẬN
#   New "Cơ sở" summary template: 4 new leading columns + shifted
#   existing columns with updated totals.
# ======================================================================
$ws8 = $wb.Worksheets.Item("LỢI NHUẬN")
$ws8.Columns("B:E").Insert()

$ws8.Range("A1").Value = "Cơ sở"
$ws8.Range("B1").Value = "Tổng đơn giá"
$ws8.Range("C1").Value = "Đã thanh toán"
$ws8.Range("D1").Value = "Tỉ lệ thanh toán"
$ws8.Range("E1").Value = "Tỉ lệ nợ"
$ws8.Range("F1").Value = "Thu nợ"
$ws8.Range("G1").Value = "Tổng doanh thu"
$ws8.Range("H1").Value = "Chi tiêu"
$ws8.Range("I1").Value = "Quỹ lương"
$ws8.Range("J1").Value = "Tổng chi phí"
$ws8.Range("K1").Value = "Lợi nhuận"
$ws8.Range("L1").Value = "Tỉ lệ lợi nhuận"

$ws8.Range("A2").Value = "CẦN THƠ"
$ws8.Range("B2").Value = 22550000
$ws8.Range("C2").Value = 22550000
$ws8.Range("D2").Value = 1
$ws8.Range("E2").Value = 0
$ws8.Range("F2").Value = 7500000
$ws8.Range("G2").Value = 30050000
$ws8.Range("H2").Value = 5296500
$ws8.Range("I2").Value = 16225023.80952381
$ws8.Range("J2").Value = 21521523.80952381
$ws8.Range("K2").Value = 8528476.19047619
$ws8.Range("L2").Value = 0.2838095238095238
